$d = $word.ActiveDocument

# The document currently has the "_GoBack" bookmark (empty, collapsed) in
# the third paragraph, and a page-break run in the first paragraph. We need
# to effectively move the bookmark to the first paragraph (replacing the
# page break) and leave the third paragraph completely empty.

# 1. Remove the existing "_GoBack" bookmark from paragraph 3.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Remove the page-break character from paragraph 1, keeping its
#    (now empty) paragraph mark intact.
$p1 = $d.Paragraphs.Item(1)
$breakRange = $d.Range($p1.Range.Start, $p1.Range.End - 1)
if ($breakRange.Start -ne $breakRange.End) {
    $breakRange.Delete()
}

# 3. Re-create the "_GoBack" bookmark, now collapsed at the start of the
#    (now empty) first paragraph.
$p1 = $d.Paragraphs.Item(1)
$bmRange = $d.Range($p1.Range.Start, $p1.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
